$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: delete the old secondary header row (row 2). This shifts data
# rows 3..10 up to 2..9 and merges the two header rows into one (row 1 keeps
# its own cells; we'll overwrite it completely below).
$ws.Rows(2).Delete()

# --- Step 2: rewrite the header row (row 1) completely with the new
# 11-column header scheme.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 keep the plain default style (no explicit format), matching the
# original workbook default.
$ws.Range("A1:E1").ClearFormats()

# F1:K1 get the "General number format, 9pt Arial font" look used
# elsewhere in the sheet, applied through a transient named style so the
# resulting cellXf keeps xfId=0 (Normal) - mirrors how the authoring tool
# produced the workbook.
$tmpStyle = $wb.Styles.Add("TmpHeaderStyle")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TmpHeaderStyle"
$tmpStyle.Delete()

# --- Step 3: the selection Excel left behind after the edit.
$ws.Range("A2:K2").Select()
